$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 137, shifting existing rows 137-139 down to 138-140.
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new weekly record.
$ws.Cells.Item(137, 1).Value = 7
$ws.Cells.Item(137, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(137, 3).Value = "Ñuble"
$ws.Cells.Item(137, 4).Value = 44448
$ws.Cells.Item(137, 5).Value = 16
$ws.Cells.Item(137, 6).Value = 100112043
$ws.Cells.Item(137, 7).Value = "Pepino ensalada"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 160
$ws.Cells.Item(137, 11).Value = 17000
$ws.Cells.Item(137, 12).Value = 18000
$ws.Cells.Item(137, 13).Value = 17500
$ws.Cells.Item(137, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(137, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(137, 16).Value = 292
$ws.Cells.Item(137, 17).Value = 60
$ws.Cells.Item(137, 18).Value = "Hortaliza"

# Match the date-time number format used by the rest of column D.
$ws.Cells.Item(137, 4).NumberFormat = $ws.Cells.Item(138, 4).NumberFormat
